# Apply corrected weights & balance figures across several sheets of the
# Weights.xlsx workbook (errors fixed in the weights and balance analyses).

$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS -------------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value  = 24956.15829483744
$ws.Range("C7").Value  = 24560.15829483744
$ws.Range("C8").Value  = 22460.542465353697
$ws.Range("C12").Value = 20001.666514896573
$ws.Range("C13").Value = 19605.666514896573
$ws.Range("C14").Value = 12873.666514896577
$ws.Range("C15").Value = 11644.122869396575
$ws.Range("C16").Value = 12058.166869396577
$ws.Range("C20").Value = 244736.30974206753
$ws.Range("C21").Value = 240852.87634206755
$ws.Range("C22").Value = 220262.67876786078
$ws.Range("C26").Value = 196149.34292831045
$ws.Range("C27").Value = 192265.90952831044
$ws.Range("C28").Value = 126247.54172831048
$ws.Range("C29").Value = 114189.83753716788
$ws.Range("C30").Value = 118250.22212976794

# --- FUSELAGE ---------------------------------------------------------
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C7").Value  = 2371.0
$ws.Range("D7").Value  = -20.966666666666665
$ws.Range("C8").Value  = 3023.0
$ws.Range("D8").Value  = 0.7666666666666667
$ws.Range("C9").Value  = 2676.0
$ws.Range("D9").Value  = -10.8
$ws.Range("C12").Value = 2915.5
$ws.Range("D12").Value = -2.816666666666666

# --- WING ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C7").Value  = 2407.0
$ws.Range("D7").Value  = 20.35
$ws.Range("C8").Value  = 1792.0
$ws.Range("D8").Value  = -10.4
$ws.Range("C9").Value  = 2027.0
$ws.Range("D9").Value  = 1.35
$ws.Range("C12").Value = 2233.0
$ws.Range("D12").Value = 11.65
$ws.Range("C13").Value = 1892.428571428571
$ws.Range("D13").Value = -5.3785714285714326

# --- HORIZONTAL TAIL ------------------------------------------------------
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C9").Value  = 133.0
$ws.Range("D9").Value  = -77.83333333333333
$ws.Range("C10").Value = 180.0
$ws.Range("D10").Value = -69.99999999999999

# --- VERTICAL TAIL --------------------------------------------------------
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C8").Value = 228.0
$ws.Range("D8").Value = -43.0
$ws.Range("C9").Value = 279.0
$ws.Range("D9").Value = -30.25

# --- LANDING GEARS --------------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 820.0
$ws.Range("D5").Value = 64.0
$ws.Range("C6").Value = 1002.0
$ws.Range("D6").Value = 100.4
$ws.Range("C7").Value = 1134.0
$ws.Range("D7").Value = 126.8
$ws.Range("C8").Value = 987.0
$ws.Range("D8").Value = 97.4
$ws.Range("C9").Value = 985.75
$ws.Range("D9").Value = 97.14999999999998
